# Auto-generated Excel COM-interop script to apply the Moogle_Profits market-data refresh diff.
# Updates numeric price/profit columns (H-N) for specific leve rows across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2068.913
$ws.Range("J17").Value = 2068.913
$ws.Range("L17").Value = 6206.739
$ws.Range("N17").Value = -6542.739
# Row 18
$ws.Range("H18").Value = 2000441
$ws.Range("I18").Value = 514.5
$ws.Range("J18").Value = 14000000
$ws.Range("K18").Value = 514.5
$ws.Range("L18").Value = 14000000
$ws.Range("M18").Value = -230.5
$ws.Range("N18").Value = -14000568
# Row 43
$ws.Range("H43").Value = 6018.0356
$ws.Range("J43").Value = 3800.6
$ws.Range("L43").Value = 3800.6
$ws.Range("N43").Value = -3938.6
# Row 100
$ws.Range("H100").Value = 2146.3333
$ws.Range("I100").Value = 720
$ws.Range("J100").Value = 4999
$ws.Range("K100").Value = 720
$ws.Range("L100").Value = 4999
$ws.Range("M100").Value = -179
$ws.Range("N100").Value = -6081
# Row 111
$ws.Range("H111").Value = 155410.72
$ws.Range("I111").Value = 2750
$ws.Range("J111").Value = 216475
$ws.Range("K111").Value = 8250
$ws.Range("L111").Value = 649425
$ws.Range("M111").Value = -5183
$ws.Range("N111").Value = -655559
# Row 127
$ws.Range("H127").Value = 1255.375
$ws.Range("I127").Value = 507.16666
$ws.Range("K127").Value = 1521.49998
$ws.Range("M127").Value = 3438.50002
# Row 131
$ws.Range("H131").Value = 3094.5
$ws.Range("I131").Value = 3094.5
$ws.Range("K131").Value = 9283.5
$ws.Range("M131").Value = -4243.5
# Row 132
$ws.Range("H132").Value = 2092.1636
$ws.Range("I132").Value = 1924.925
$ws.Range("J132").Value = 2538.1333
$ws.Range("K132").Value = 5774.775
$ws.Range("L132").Value = 7614.3999
$ws.Range("M132").Value = -3244.775
$ws.Range("N132").Value = -12674.3999
# Row 138
$ws.Range("H138").Value = 7427.9375
$ws.Range("I138").Value = 3153.6667
$ws.Range("J138").Value = 9992.5
$ws.Range("K138").Value = 9461.000100000001
$ws.Range("L138").Value = 29977.5
$ws.Range("M138").Value = -4321.000100000001
$ws.Range("N138").Value = -40257.5
# Row 141
$ws.Range("H141").Value = 3713.4285
$ws.Range("I141").Value = 1371.8
$ws.Range("J141").Value = 9567.5
$ws.Range("K141").Value = 4115.4
$ws.Range("L141").Value = 28702.5
$ws.Range("M141").Value = 1064.6
$ws.Range("N141").Value = -39062.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 7097.718
$ws.Range("I61").Value = 8212.076999999999
$ws.Range("K61").Value = 8212.076999999999
$ws.Range("M61").Value = -8000.076999999999
# Row 88
$ws.Range("H88").Value = 892.2759
$ws.Range("I88").Value = 562.7273
$ws.Range("K88").Value = 562.7273
$ws.Range("M88").Value = -156.7273
# Row 91
$ws.Range("H91").Value = 892.2759
$ws.Range("I91").Value = 562.7273
$ws.Range("K91").Value = 562.7273
$ws.Range("M91").Value = 841.2727
# Row 136
$ws.Range("H136").Value = 7097.718
$ws.Range("I136").Value = 8212.076999999999
$ws.Range("K136").Value = 24636.231
$ws.Range("M136").Value = -22086.231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2411.56
$ws.Range("I20").Value = 2185.6
$ws.Range("J20").Value = 2750.5
$ws.Range("K20").Value = 2185.6
$ws.Range("L20").Value = 2750.5
$ws.Range("M20").Value = -1938.6
$ws.Range("N20").Value = -3244.5
# Row 94
$ws.Range("H94").Value = 911.6
$ws.Range("I94").Value = 889.4375
$ws.Range("K94").Value = 889.4375
$ws.Range("M94").Value = -438.4375
# Row 105
$ws.Range("H105").Value = 1179517.9
$ws.Range("I105").Value = 1484137.4
$ws.Range("J105").Value = 4557.143
$ws.Range("K105").Value = 1484137.4
$ws.Range("L105").Value = 4557.143
$ws.Range("M105").Value = -1482390.4
$ws.Range("N105").Value = -8051.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 14984.333
$ws.Range("I31").Value = 5093.6
$ws.Range("J31").Value = 22049.143
$ws.Range("K31").Value = 5093.6
$ws.Range("L31").Value = 22049.143
$ws.Range("M31").Value = -4798.6
$ws.Range("N31").Value = -22639.143
# Row 34
$ws.Range("H34").Value = 14984.333
$ws.Range("I34").Value = 5093.6
$ws.Range("J34").Value = 22049.143
$ws.Range("K34").Value = 5093.6
$ws.Range("L34").Value = 22049.143
$ws.Range("M34").Value = -4891.6
$ws.Range("N34").Value = -22453.143
# Row 41
$ws.Range("H41").Value = 33999.5
$ws.Range("I41").Value = 8000
$ws.Range("J41").Value = 59999
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 59999
$ws.Range("M41").Value = -7572
$ws.Range("N41").Value = -60855
# Row 60
$ws.Range("H60").Value = 25579.8
$ws.Range("J60").Value = 47999.5
$ws.Range("L60").Value = 47999.5
$ws.Range("N60").Value = -49021.5
# Row 68
$ws.Range("H68").Value = 99999.664
$ws.Range("J68").Value = 99999.664
$ws.Range("L68").Value = 99999.664
$ws.Range("N68").Value = -101497.664
# Row 71
$ws.Range("H71").Value = 99999.664
$ws.Range("J71").Value = 99999.664
$ws.Range("L71").Value = 299998.992
$ws.Range("N71").Value = -307486.992

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 3249.75
$ws.Range("J68").Value = 3249.75
$ws.Range("L68").Value = 9749.25
$ws.Range("N68").Value = -11371.25
# Row 71
$ws.Range("H71").Value = 3249.75
$ws.Range("J71").Value = 3249.75
$ws.Range("L71").Value = 29247.75
$ws.Range("N71").Value = -37359.75
# Row 129
$ws.Range("H129").Value = 13892315
$ws.Range("J129").Value = 18521790
$ws.Range("L129").Value = 55565370
$ws.Range("N129").Value = -55575370
# Row 130
$ws.Range("H130").Value = 7500
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
# Row 131
$ws.Range("H131").Value = 1215661.8
$ws.Range("I131").Value = 1161.1428
$ws.Range("J131").Value = 2278349.8
$ws.Range("K131").Value = 3483.4284
$ws.Range("L131").Value = 6835049.399999999
$ws.Range("M131").Value = 1556.5716
$ws.Range("N131").Value = -6845129.399999999
# Row 139
$ws.Range("H139").Value = 2222.4443
$ws.Range("I139").Value = 2060.3333
$ws.Range("J139").Value = 3033
$ws.Range("K139").Value = 6180.999899999999
$ws.Range("L139").Value = 9099
$ws.Range("M139").Value = -1040.999899999999
$ws.Range("N139").Value = -19379
# Row 140
$ws.Range("H140").Value = 2056.1025
$ws.Range("I140").Value = 2007.3334
$ws.Range("J140").Value = 2077.7778
$ws.Range("K140").Value = 6022.0002
$ws.Range("L140").Value = 6233.3334
$ws.Range("M140").Value = -842.0002000000004
$ws.Range("N140").Value = -16593.3334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 11716.5
$ws.Range("I80").Value = 2719.8
$ws.Range("K80").Value = 2719.8
$ws.Range("M80").Value = -1721.8
# Row 83
$ws.Range("H83").Value = 11716.5
$ws.Range("I83").Value = 2719.8
$ws.Range("K83").Value = 13599
$ws.Range("M83").Value = -8607
# Row 126
$ws.Range("H126").Value = 4924.6113
$ws.Range("I126").Value = 6700.5
$ws.Range("J126").Value = 4036.6667
$ws.Range("K126").Value = 20101.5
$ws.Range("L126").Value = 12110.0001
$ws.Range("M126").Value = -17631.5
$ws.Range("N126").Value = -17050.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 263101
$ws.Range("I7").Value = 263101
$ws.Range("K7").Value = 263101
$ws.Range("M7").Value = -262989
# Row 55
$ws.Range("H55").Value = 1125.4706
$ws.Range("I55").Value = 252.33333
$ws.Range("K55").Value = 252.33333
$ws.Range("M55").Value = -79.33332999999999
# Row 126
$ws.Range("H126").Value = 263101
$ws.Range("I126").Value = 263101
$ws.Range("K126").Value = 789303
$ws.Range("M126").Value = -786833
# Row 132
$ws.Range("H132").Value = 2299.5757
$ws.Range("I132").Value = 1343.963
$ws.Range("K132").Value = 4031.889
$ws.Range("M132").Value = -1501.889

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 7994.5
$ws.Range("I52").Value = 5433.3335
$ws.Range("J52").Value = 15678
$ws.Range("K52").Value = 5433.3335
$ws.Range("L52").Value = 15678
$ws.Range("M52").Value = -5207.3335
$ws.Range("N52").Value = -16130
# Row 122
$ws.Range("H122").Value = 2504.3215
$ws.Range("I122").Value = 2486
$ws.Range("K122").Value = 7458
$ws.Range("M122").Value = -5008
